# Fix bug for upgrade building: correct the typo "dragonEyire" -> "dragonEyrie"
# and restore the active selection to the edited cell (B5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "dragonEyrie"
$ws.Range("B5").Select()
